$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the format (incl. the mmm-yy date style) of the last existing
# date cell down onto the three new rows before writing their values,
# so the new dates pick up the same style index as the rest of column A.
$ws.Range("A106").Copy()
$ws.Range("A107:A109").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New monthly seasonality data points (Oct/Nov/Dec 2022).
$ws.Range("A107").Value = 44835
$ws.Range("B107").Value = 646
$ws.Range("A108").Value = 44866
$ws.Range("B108").Value = 613
$ws.Range("A109").Value = 44896
$ws.Range("B109").Value = 581

# Matches the author's recorded selection after entering the new rows.
$ws.Range("C107").Select()
